$d = $word.ActiveDocument

# Locate the end of the "checkout view template." paragraph (last bullet
# of the original list) and insert a fresh paragraph mark right after it.
$rng = $d.Content
$rng.Find.Execute("Added Input validation error message displaying to the checkout view template.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertParagraphAfter()

# Re-acquire the freshly created (empty) paragraph by its index so the
# Range we act on is not stale.
$insertIndex = $rng.Paragraphs(1).Index + 1
$target = $d.Paragraphs($insertIndex).Range

$pkg = @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr>
<w:ilvl w:val="0"/>
<w:numId w:val="1"/>
</w:numPr>
</w:pPr>
<w:r><w:t xml:space="preserve">Added authentication and authorization to the project. Added Nuget Packages - Microsoft.AspNetCore.Identity.EntityFrameworkCore and </w:t></w:r>
<w:r><w:t>Microsoft.AspNetCore.Identity.</w:t></w:r>
<w:r><w:t xml:space="preserve">UI. Updated the </w:t></w:r>
<w:r><w:t>AppDbContext.cs</w:t></w:r>
<w:r><w:t xml:space="preserve"> class by changing the inherit</w:t></w:r>
<w:r><w:t>ance</w:t></w:r>
<w:r><w:t xml:space="preserve"> the AppDbContext from IdentityDbContext which is the bridge between the application and the database.</w:t></w:r>
</w:p>
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr>
<w:ilvl w:val="0"/>
<w:numId w:val="1"/>
</w:numPr>
</w:pPr>
<w:r><w:t>Added migration to add identity to the database and updated it. Updated the startup class by adding the use authentication middleware.</w:t></w:r>
</w:p>
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr>
<w:ilvl w:val="0"/>
<w:numId w:val="1"/>
</w:numPr>
</w:pPr>
<w:r><w:lastRenderedPageBreak/><w:t>Added new scaffold identity item with account login, logout and register.</w:t></w:r>
<w:r><w:t xml:space="preserve"> Added configure services to configure identity service in the startup.cs class.</w:t></w:r>
<w:r><w:t xml:space="preserve"> Added a link to the login page using _Login partial view in the _Layout file.</w:t></w:r>
<w:r><w:t xml:space="preserve"> Enable authorization</w:t></w:r>
<w:r><w:t xml:space="preserve"> by adding [Authorize] attribute in the properties.</w:t></w:r>
<w:r><w:t xml:space="preserve"> Edited the Order controller with authorization. Configured authorization in the startup.cs file.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$target.InsertXML($pkg)
